$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Inventory
$ws.Range("B4").Value = 1305000000
$ws.Range("C4").Value = 1338000000
$ws.Range("D4").Value = 1343000000
$ws.Range("E4").Value = 1402000000
$ws.Range("F4").Value = 1346000000

# Row 12: Accounts Payable
$ws.Range("B12").Value = 215000000
$ws.Range("C12").Value = 207000000
$ws.Range("D12").Value = 210000000
$ws.Range("E12").Value = 194000000
$ws.Range("F12").Value = 212000000

# Row 19: Long Term Tax Liability (Deferred)
$ws.Range("B19").Value = 103000000
$ws.Range("C19").Value = 102000000
$ws.Range("D19").Value = 103000000
$ws.Range("E19").Value = 100000000
$ws.Range("F19").Value = 100000000
